$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6, column A: task id #100004 (reuses freed shared-string slot
# left behind once B5's old text is no longer referenced)
$ws.Range("A6").Value = "#100004"

# Row 5, column B: update task description for #100003
$ws.Range("B5").Value = "Add crud features for Trucks"

# New row 6, column B: task description for #100004
$ws.Range("B6").Value = "Add crud features for Mop and extendedMop"

# Move the active selection down to B7, matching the post-edit cursor position
$ws.Range("B7").Select()
